# AutoCommit_14 декабря 2023 г. 12:29:42_SibNout2023
# Fill in previously-missing homework grades (value 5) for several
# students, copying the cell format from an existing "filled" grade cell
# in the same style family so the thick-border formatting matches, then
# move the active selection to Q5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Fill-Grade($addr, $formatSource, $value) {
    $ws.Range($formatSource).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
    $ws.Range($addr).Value = $value
}

# Row 5: L5 was an empty "style 6" placeholder -> becomes a filled grade
# cell with the thick-left-border style (style 8, e.g. like J10/K10/L10).
Fill-Grade "L5" "J10" 5

# Row 14: J14 newly appears with the thick-left-border style (style 8);
# L14 moves from the empty placeholder style to the borderless
# centered style (style 10, e.g. like K12/L12).
Fill-Grade "J14" "J10" 5
Fill-Grade "L14" "K12" 5

# Row 15: extra column T gets a plain value-only cell (no special style).
$ws.Range("T15").Value = 5

# Row 20: K20 newly appears and L20 is filled in, both using the
# borderless centered style (style 10).
Fill-Grade "K20" "K12" 5
Fill-Grade "L20" "K12" 5

# Row 21: L21 filled in with the borderless centered style (style 10).
Fill-Grade "L21" "K12" 5

# Row 23: extra column T gets a plain value-only cell.
$ws.Range("T23").Value = 5

# Row 29: K29, L29, M29 all filled in with the thick-left-border style
# (style 8), and column T gets a plain value-only cell.
Fill-Grade "K29" "J10" 5
Fill-Grade "L29" "J10" 5
Fill-Grade "M29" "J10" 5
$ws.Range("T29").Value = 5

$excel.CutCopyMode = 0

# Move the active selection in the frozen bottom-right pane to Q5.
$ws.Range("Q5").Select()
